$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 5 label, matching the style used by the other year cells in column A
$ws.Range("A5").Value = "2021年"
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)

# New row 5 numeric data (columns B:AT)
$values = @(31.708,42.831,49.096,32.737,40,26.963,32.578,35.019,30.614,43.645,30.904,29.518,26.234,28.616,27.292,30.125,21.398,18.46,31.101,19.12,37.302,46.825,17.419,23.772,19.603,18.602,39.837,20.48,29.198,32.867,24.827,19.683,46.019,38.129,27.513,30.462,18.526,29.31,30.796,37.106,25.659,17.939,33.409,26.541,17.503)

for ($i = 0; $i -lt $values.Length; $i++) {
    $col = $i + 2
    $ws.Cells.Item(5, $col).Value = $values[$i]
}
